$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count())
$newSheet = $wb.Worksheets.Add([Type]::Missing, $last)
$newSheet.Name = "BC"
$newSheet.Range("A1").Value = "Topics"
$newSheet.Range("A1:B1").Merge()
Write-Host "merged"
